$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new cell text via a formula (CHAR(10) for embedded newlines), then
# convert the formula result to a plain value in place. Using PasteSpecial
# (paste-values) here -- rather than re-assigning .Value with the literal
# string -- avoids Excel's "recalculate row height for wrapped text" pass,
# matching the target row (no ht/customHeight attribute).
$ws.Range("A1").Formula = @'
="questions = ["&CHAR(10)&"    {"&CHAR(10)&"        ""title"": ""You are given an array of integers, arr and its size n. You are also given an empty result array, rarr. Write the function arrayProducts that takes n and arr as inputs. It fills in the result array so that the element at every position is the product of the elements at all the other positions.Example 1Input:n=4arr= 1 2 3 4rarr= emptyOutput:rarr = 24 12 8 6Explanation: The first element in the result array = product of the remaining elements = 2 x 3 x 4 = 24The second element = 1 x 3 x 4 = 12The third element = 1 x 2 x 4 = 8The fourth element = 1 x 2 x 3 = 6Example 2Input:n=5arr= -1 1 3 5 7rarr= emptyOutput:rarr = 105 -105 -35 -21 -15Explanation: The elements will be [1x3x5x7, -1x3x5x7, -1x1x5x7, -1x1x3x7, -1x1x3x5] = [105, -105, -35, -21, -15]"","&CHAR(10)&"        ""ques_type"": null,"&CHAR(10)&"        ""options"": [],"&CHAR(10)&"        ""score"": null"&CHAR(10)&"    }"&CHAR(10)&"]"
'@
$ws.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Drop the old second row/cell entirely.
$ws.Range("A2").ClearContents()

# A1 no longer needs the bold/bordered/centered style that used to mark the
# header cell -- put it back to the workbook's default "Normal" style.
$ws.Range("A1").Style = "Normal"
